$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-typed team name "Olbia-virtusverona" -> "Olbia-virtus"
$ws.Range("A2").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "Olbia-virtus"

# Add the new "Piede" column (replacing what used to be a separate list)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "Piede"

for ($r = 2; $r -le 16; $r++) {
    $ws.Range("H1").Copy($ws.Cells.Item($r, 9))
    $ws.Cells.Item($r, 9).Value = "Destro"
}
